# Auto-generated edit script
# Applies numeric corrections to the Leve profit-tracking columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, LTW, and WVR sheets, per the
# scheduled pricing-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6333.3335
$ws.Range("J51").Value = 6500
$ws.Range("L51").Value = 6500
$ws.Range("N51").Value = -7468
$ws.Range("H100").Value = 2121.75
$ws.Range("I100").Value = 912.2857
$ws.Range("J100").Value = 3815
$ws.Range("K100").Value = 912.2857
$ws.Range("L100").Value = 3815
$ws.Range("M100").Value = -371.2857
$ws.Range("N100").Value = -4897
$ws.Range("H112").Value = 1813.3846
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1813.3846
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 5440.1538
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -7656.1538
$ws.Range("H135").Value = 52632540
$ws.Range("I135").Value = 1074.2858
$ws.Range("J135").Value = 200000640
$ws.Range("K135").Value = 9668.572200000001
$ws.Range("L135").Value = 1800005760
$ws.Range("M135").Value = -7133.572200000001
$ws.Range("N135").Value = -1800010830
$ws.Range("H138").Value = 2047.3
$ws.Range("I138").Value = 1830.5135
$ws.Range("J138").Value = 2396.0435
$ws.Range("K138").Value = 5491.5405
$ws.Range("L138").Value = 7188.130500000001
$ws.Range("M138").Value = -351.5405000000001
$ws.Range("N138").Value = -17468.1305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 398132.94
$ws.Range("I2").Value = 618326.1
$ws.Range("J2").Value = 1785.2
$ws.Range("K2").Value = 618326.1
$ws.Range("L2").Value = 1785.2
$ws.Range("M2").Value = -618213.1
$ws.Range("N2").Value = -2011.2
$ws.Range("H32").Value = 4632.3
$ws.Range("I32").Value = 3747.0698
$ws.Range("J32").Value = 10070.143
$ws.Range("K32").Value = 3747.0698
$ws.Range("L32").Value = 10070.143
$ws.Range("M32").Value = -3460.0698
$ws.Range("N32").Value = -10644.143
$ws.Range("H116").Value = 398132.94
$ws.Range("I116").Value = 618326.1
$ws.Range("J116").Value = 1785.2
$ws.Range("K116").Value = 618326.1
$ws.Range("L116").Value = 1785.2
$ws.Range("M116").Value = -616032.1
$ws.Range("N116").Value = -6373.2
$ws.Range("H132").Value = 1359.5714
$ws.Range("I132").Value = 1386.1177
$ws.Range("K132").Value = 4158.3531
$ws.Range("M132").Value = -1628.3531

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 398132.94
$ws.Range("I3").Value = 618326.1
$ws.Range("J3").Value = 1785.2
$ws.Range("K3").Value = 618326.1
$ws.Range("L3").Value = 1785.2
$ws.Range("M3").Value = -618212.1
$ws.Range("N3").Value = -2013.2
$ws.Range("H94").Value = 728.3077
$ws.Range("I94").Value = 755.63635
$ws.Range("J94").Value = 578
$ws.Range("K94").Value = 755.63635
$ws.Range("L94").Value = 578
$ws.Range("M94").Value = -304.63635
$ws.Range("N94").Value = -1480
$ws.Range("H134").Value = 6623.0415
$ws.Range("I134").Value = 7802.8423
$ws.Range("K134").Value = 23408.5269
$ws.Range("M134").Value = -20873.5269

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3131.6428
$ws.Range("I31").Value = 3264.2856
$ws.Range("J31").Value = 2999
$ws.Range("K31").Value = 3264.2856
$ws.Range("L31").Value = 2999
$ws.Range("M31").Value = -2969.2856
$ws.Range("N31").Value = -3589
$ws.Range("H34").Value = 3131.6428
$ws.Range("I34").Value = 3264.2856
$ws.Range("J34").Value = 2999
$ws.Range("K34").Value = 3264.2856
$ws.Range("L34").Value = 2999
$ws.Range("M34").Value = -3062.2856
$ws.Range("N34").Value = -3403
$ws.Range("H50").Value = 14140
$ws.Range("J50").Value = 14140
$ws.Range("L50").Value = 14140
$ws.Range("N50").Value = -15390
$ws.Range("H99").Value = 2768.8333
$ws.Range("I99").Value = 2106.7144
$ws.Range("K99").Value = 2106.7144
$ws.Range("M99").Value = -608.7143999999998
$ws.Range("H126").Value = 2768.8333
$ws.Range("I126").Value = 2106.7144
$ws.Range("K126").Value = 6320.1432
$ws.Range("M126").Value = -3850.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 467
$ws.Range("I5").Value = 358.91666
$ws.Range("J5").Value = 899.3333
$ws.Range("K5").Value = 1076.74998
$ws.Range("L5").Value = 2697.9999
$ws.Range("M5").Value = -964.7499800000001
$ws.Range("N5").Value = -2921.9999
$ws.Range("H11").Value = 1340
$ws.Range("I11").Value = 1200
$ws.Range("J11").Value = 1480
$ws.Range("K11").Value = 3600
$ws.Range("L11").Value = 4440
$ws.Range("M11").Value = -3460
$ws.Range("N11").Value = -4720
$ws.Range("H33").Value = 85
$ws.Range("I33").Value = 59.333332
$ws.Range("J33").Value = 162
$ws.Range("K33").Value = 355.999992
$ws.Range("L33").Value = 972
$ws.Range("M33").Value = -72.99999200000002
$ws.Range("N33").Value = -1538
$ws.Range("H122").Value = 1551.6666
$ws.Range("J122").Value = 1551.6666
$ws.Range("L122").Value = 13964.9994
$ws.Range("N122").Value = -18864.9994
$ws.Range("H131").Value = 65080.273
$ws.Range("J131").Value = 101854.71
$ws.Range("L131").Value = 305564.13
$ws.Range("N131").Value = -315644.13
$ws.Range("H135").Value = 467
$ws.Range("I135").Value = 358.91666
$ws.Range("J135").Value = 899.3333
$ws.Range("K135").Value = 3230.24994
$ws.Range("L135").Value = 8093.9997
$ws.Range("M135").Value = -695.2499399999997
$ws.Range("N135").Value = -13163.9997
$ws.Range("H137").Value = 4537.0835
$ws.Range("J137").Value = 6834
$ws.Range("L137").Value = 20502
$ws.Range("N137").Value = -30702
$ws.Range("H139").Value = 1990.3529
$ws.Range("I139").Value = 1998
$ws.Range("K139").Value = 5994
$ws.Range("M139").Value = -854
$ws.Range("H140").Value = 2459.3684
$ws.Range("I140").Value = 1548.6364
$ws.Range("K140").Value = 4645.9092
$ws.Range("M140").Value = 534.0907999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2138.4348
$ws.Range("I22").Value = 2320.2856
$ws.Range("J22").Value = 1855.5555
$ws.Range("K22").Value = 2320.2856
$ws.Range("L22").Value = 1855.5555
$ws.Range("M22").Value = -2025.2856
$ws.Range("N22").Value = -2445.5555
$ws.Range("H27").Value = 2138.4348
$ws.Range("I27").Value = 2320.2856
$ws.Range("J27").Value = 1855.5555
$ws.Range("K27").Value = 2320.2856
$ws.Range("L27").Value = 1855.5555
$ws.Range("M27").Value = -2213.2856
$ws.Range("N27").Value = -2069.5555
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H139").Value = 42998.75
$ws.Range("J139").Value = 42998.75
$ws.Range("L139").Value = 42998.75
$ws.Range("N139").Value = -53278.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H126").Value = 7647.6113
$ws.Range("I126").Value = 8365.200000000001
$ws.Range("J126").Value = 4059.6667
$ws.Range("K126").Value = 25095.6
$ws.Range("L126").Value = 12179.0001
$ws.Range("M126").Value = -22625.6
$ws.Range("N126").Value = -17119.0001
$ws.Range("H139").Value = 72586.2
$ws.Range("J139").Value = 72586.2
$ws.Range("L139").Value = 72586.2
$ws.Range("N139").Value = -82866.2
